# Add an Excel "Data Model"-style header row and shift the existing
# helper rows/columns to make room for it, per commit:
# "Add Excel Data Model and update Helpers"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current data, then insert two new
# blank columns to the left of the current data, pushing the existing
# Ticket/Description/Priority table from A1:C4 down to C2:E5.
$null = $ws.Rows.Item(1).Insert()
$null = $ws.Columns.Item(1).Insert()
$null = $ws.Columns.Item(1).Insert()

# Populate the new header row (row 1) with the Data Model column names.
$ws.Range("A1").Value = "Board"
$ws.Range("B1").Value = "List"
$ws.Range("D1").Value = "Description"
$ws.Range("C1").Value = "Tittle "
$ws.Range("E1").Value = "Label"

# Match the saved selection/active cell.
$null = $ws.Range("D2").Select()
